$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "자주 묻는 질문"

$ws.Range("A2").Value = "_질문"
$ws.Range("B2").Value = "_답변"
$ws.Range("C2").Value = "_소스_uri"

$ws.Range("A3").Value = "자동차 보험료를 결정하는 요인은 무엇인가요?"
$ws.Range("B3").Value = "보험료는 운전 기록, 연령, 위치, 차량 유형, 보장 옵션 등 다양한 요인의 영향을 받습니다. 보험료를 결정하기 위해 맞춤형 견적을 만들어 드릴 수 있습니다."
$ws.Range("C3").Value = "https://www.anycompany.com/premuims"

$ws.Range("A4").Value = "어떤 유형의 보험을 제공하나요?"
$ws.Range("B4").Value = "당사는 책임, 충돌, 종합 보험, 보험에 가입하지 않았거나 보험이 부족한 운전자, 의료비 지급 등 다양한 보장 옵션을 제공합니다."
$ws.Range("C4").Value = "https://www.anycompany.com/coverage"

$ws.Range("A5").Value = "자동차 보험료를 낮추려면 어떻게 해야 하나요?"
$ws.Range("B5").Value = "운전 기록을 잘 유지하고, 방어적인 운전 코스를 수강하고, 번들링 정책을 적용하고, 할인을 받는 등 여러 가지 방법으로 요금을 낮출 수 있습니다. 구체적인 요금에 대해서는 상담원에게 문의할 수 있습니다."
$ws.Range("C5").Value = "https://www.anycompany.com/rates"

$ws.Range("A6").Value = "보험계약에 새 드라이버를 추가할 수 있나요?"
$ws.Range("B6").Value = "예. 보험계약에 새 드라이버를 추가할 수 있습니다. 보험료는 신규 운전자의 연령, 운전 기록 및 기타 요인에 따라 달라질 수 있습니다."
$ws.Range("C6").Value = "https://www.anycompany.com/policy"

$ws.Range("A7").Value = "사고 후에는 어떻게 해야 하나요?"
$ws.Range("B7").Value = "먼저 모든 사람의 안전을 확인한 다음 관련된 다른 당사자와 정보를 교환하세요. 가능하면 사진을 찍고 경찰에 신고하세요. 즉시 당사에 연락하여 사고를 신고하고 청구 절차를 시작하십시오."
$ws.Range("C7").Value = "https://www.anycompany.com/claims"

$ws.Range("A8").Value = "공제액이란 무엇이며 어떻게 적용되나요?"
$ws.Range("B8").Value = "공제액은 보험이 적용되기 전에 지불해야 하는 금액입니다. 예를 들어, 공제액이 500 달러이고 수리비가 2000 USD인 경우 500을 지불하고 나머지 1,500을 보험에서 부담합니다."
$ws.Range("C8").Value = "https://www.anycompany.com/deductibles"

$ws.Range("A7").Select()
